# Hortaliza, Macroferia Regional de Talca - Repollo
# A new weekly price-report row is inserted at row 470 (pushing the existing
# rows 470-553 down to 471-554). The new row carries the same
# Mercado/Region/Producto/Variedad/Calidad/Unidad/Origen/Clasificacion
# metadata as the row immediately below it (the row that used to be 470),
# but with its own Fecha/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 470; everything from the old row 470
# downward shifts to row 471 onward.
$ws.Rows.Item(470).Insert()

# Seed the new row 470 with the same values as the row below it (which is
# the former row 470), then overwrite the date/volume/price columns with
# the new record's values.
$src = $ws.Range("A471:R471")
$dst = $ws.Range("A470:R470")
$src.Copy($dst)

$ws.Range("D470").Value = 45180
$ws.Range("J470").Value = 3000
$ws.Range("K470").Value = 600
$ws.Range("L470").Value = 600
$ws.Range("M470").Value = 600
$ws.Range("P470").Value = 600
